# Apply the updated crypto price/volume figures captured by the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.865.87'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.53'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.88'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.66%  '
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0609'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.874.51'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.631.82'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.590'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.51'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.96%  '
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.866.47'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.61'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.88'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.82%  '
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.59'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.53'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.63'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0494'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +1.85%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.423.27'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  +3.77%  '
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('E37').Value = '  -5.64%  '
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.88'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +11.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.562'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.835'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0500'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.781.51'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.94'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -9.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '93.65'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.47%  '
$ws.Range('E51').Value = '  +0.13%  '
